$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a "Price" (column D) cell to a text string, forcing text
# storage (NumberFormat "@") for values that Excel would otherwise
# auto-convert to a number (single decimal point, no thousands dots).
function Set-PriceText($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.987.37"
$ws.Range("E2").Value = "  -0.98%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.441.32"
$ws.Range("E3").Value = "  -1.45%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
Set-PriceText "D5" "575.67"
$ws.Range("E5").Value = "  -1.25%  "

# Row 6 - Solana
Set-PriceText "D6" "159.60"
$ws.Range("E6").Value = "  -1.08%  "

# Row 8 - XRP
Set-PriceText "D8" "0.590"
$ws.Range("E8").Value = "  -2.93%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "3.442.83"
$ws.Range("E9").Value = "  -1.63%  "

# Row 10 - Toncoin
Set-PriceText "D10" "7.25"
$ws.Range("E10").Value = "  -0.94%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -2.52%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.10%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.031.98"
$ws.Range("E13").Value = "  -1.62%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.45%  "

# Row 15 - ShibaInu
Set-PriceText "D15" "0.0000189"
$ws.Range("E15").Value = "  -3.31%  "

# Row 16 - Avalanche
Set-PriceText "D16" "27.90"
$ws.Range("E16").Value = "  -3.06%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "64.949.27"
$ws.Range("E17").Value = "  -1.02%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.435.52"
$ws.Range("E18").Value = "  -1.28%  "

# Row 19 - Polkadot
Set-PriceText "D19" "6.38"
$ws.Range("E19").Value = "  -1.63%  "

# Row 20 - Chainlink
Set-PriceText "D20" "13.94"
$ws.Range("E20").Value = "  -2.83%  "

# Row 21 - BitcoinCash
Set-PriceText "D21" "382.99"
$ws.Range("E21").Value = "  -1.58%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -3.42%  "

# Row 23 - Polygon
Set-PriceText "D23" "0.551"
$ws.Range("E23").Value = "  -0.90%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.21%  "

# Row 25 - Litecoin
Set-PriceText "D25" "72.16"

# Row 26 - PEPE
Set-PriceText "D26" "0.0000120"
$ws.Range("E26").Value = "  -4.00%  "

# Row 27 - InternetComputer(DFINITY)
Set-PriceText "D27" "9.92"
$ws.Range("E27").Value = "  -3.40%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -0.76%  "

# Row 29 - Binance-PegBSC-USD
Set-PriceText "D29" "0.999"
$ws.Range("E29").Value = "  +0.10%  "

# Row 30 - Fetch.AI
Set-PriceText "D30" "1.49"
$ws.Range("E30").Value = "  +1.52%  "

# Row 31 - NEARProtocol
$ws.Range("E31").Value = "  -2.49%  "

# Row 32 - PancakeSwap
Set-PriceText "D32" "2.02"
$ws.Range("E32").Value = "  -2.44%  "

# Row 33 - EthereumClassic
Set-PriceText "D33" "23.34"
$ws.Range("E33").Value = "  -1.70%  "

# Row 34 - Aptos
Set-PriceText "D34" "7.08"
$ws.Range("E34").Value = "  -2.52%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  -0.30%  "

# Row 36 - Monero
Set-PriceText "D36" "160.93"
$ws.Range("E36").Value = "  -1.24%  "

# Row 37 - Stacks
$ws.Range("E37").Value = "  -1.82%  "

# Row 38 - Maker
$ws.Range("D38").Value = "2.913.67"
$ws.Range("E38").Value = "  -4.34%  "

# Row 39 - Hedera
Set-PriceText "D39" "0.0753"
$ws.Range("E39").Value = "  -3.19%  "

# Row 40 - RenderToken
Set-PriceText "D40" "6.81"

# Row 41 - EnergySwap
Set-PriceText "D41" "26.46"
$ws.Range("E41").Value = "  -3.32%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  +0.04%  "

# Row 43 - OKB
Set-PriceText "D43" "43.28"

# Row 44 - VeChain
Set-PriceText "D44" "0.0319"
$ws.Range("E44").Value = "  -2.26%  "

# Rows 45/46 swap: Mantle <-> InjectiveProtocol
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-PriceText "D45" "26.18"
$ws.Range("E45").Value = "  +0.88%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-PriceText "D46" "0.775"
$ws.Range("E46").Value = "  -0.70%  "

# Row 47 - dogwifhat
$ws.Range("E47").Value = "  +2.07%  "

# Rows 48/49 swap: Bittensor <-> ONDO
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-PriceText "D48" "1.08"
$ws.Range("E48").Value = "  -3.30%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-PriceText "D49" "316.33"
$ws.Range("E49").Value = "  -1.15%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  -3.49%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -2.95%  "
